$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.372720956802368
$ws.Range("B1").Value = 2.151985883712769
$ws.Range("C1").Value = 4.84665584564209
$ws.Range("D1").Value = 3.536446094512939
$ws.Range("E1").Value = 1.25377345085144
